$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1491935483870968
$ws.Range("C2").Value = 0.6612903225806451
$ws.Range("J2").Value = 0.02016129032258064
$ws.Range("P2").Value = 0.08870967741935484
$ws.Range("S2").Value = 0.08064516129032258
$ws.Range("B3").Value = 0.005847953216374269
$ws.Range("C3").Value = 0.02923976608187134
$ws.Range("J3").Value = 0.02339181286549707
$ws.Range("O3").Value = 0.005847953216374269
$ws.Range("P3").Value = 0.7894736842105263
$ws.Range("S3").Value = 0.1461988304093567
$ws.Range("J4").Value = 0.07317073170731707
$ws.Range("P4").Value = 0.6097560975609756
$ws.Range("S4").Value = 0.3170731707317073
$ws.Range("B6").Value = 0.0379746835443038
$ws.Range("F6").Value = 0.02953586497890295
$ws.Range("J6").Value = 0.3206751054852321
$ws.Range("O6").Value = 0.02531645569620253
$ws.Range("Q6").Value = 0.1772151898734177
$ws.Range("R6").Value = 0.06751054852320675
$ws.Range("S6").Value = 0.3417721518987342
$ws.Range("B7").Value = 0.08808290155440414
$ws.Range("D7").Value = 0.02590673575129534
$ws.Range("E7").Value = 0.005181347150259068
$ws.Range("F7").Value = 0.03626943005181347
$ws.Range("J7").Value = 0.1450777202072539
$ws.Range("O7").Value = 0.0310880829015544
$ws.Range("Q7").Value = 0.1968911917098446
$ws.Range("R7").Value = 0.08808290155440414
$ws.Range("S7").Value = 0.383419689119171
$ws.Range("B8").Value = 0.08527131782945736
$ws.Range("D8").Value = 0.01937984496124031
$ws.Range("E8").Value = 0.001937984496124031
$ws.Range("F8").Value = 0.0755813953488372
$ws.Range("J8").Value = 0.09108527131782945
$ws.Range("O8").Value = 0.01744186046511628
$ws.Range("Q8").Value = 0.1782945736434109
$ws.Range("R8").Value = 0.09689922480620156
$ws.Range("S8").Value = 0.4341085271317829
$ws.Range("B9").Value = 0.09392265193370165
$ws.Range("D9").Value = 0.01104972375690608
$ws.Range("E9").Value = 0.005524861878453038
$ws.Range("F9").Value = 0.04419889502762431
$ws.Range("J9").Value = 0.1049723756906077
$ws.Range("O9").Value = 0.01657458563535912
$ws.Range("Q9").Value = 0.1823204419889503
$ws.Range("R9").Value = 0.1325966850828729
$ws.Range("S9").Value = 0.4088397790055249
$ws.Range("B10").Value = 0.09595559080095163
$ws.Range("D10").Value = 0.01982553528945281
$ws.Range("F10").Value = 0.06978588421887391
$ws.Range("J10").Value = 0.1046788263283109
$ws.Range("O10").Value = 0.01586042823156225
$ws.Range("Q10").Value = 0.2275971451229183
$ws.Range("R10").Value = 0.1015067406819984
$ws.Range("S10").Value = 0.3647898493259318
$ws.Range("G11").Value = 0.1501831501831502
$ws.Range("J11").Value = 0.07326007326007326
$ws.Range("K11").Value = 0.1868131868131868
$ws.Range("L11").Value = 0.5677655677655677
$ws.Range("S11").Value = 0.02197802197802198
$ws.Range("G12").Value = 0.7716049382716049
$ws.Range("J12").Value = 0.1604938271604938
$ws.Range("K12").Value = 0.006172839506172839
$ws.Range("L12").Value = 0.03703703703703703
$ws.Range("S12").Value = 0.02469135802469136
$ws.Range("G13").Value = 0.775
$ws.Range("J13").Value = 0.175
$ws.Range("S13").Value = 0.05
$ws.Range("F15").Value = 0.03139013452914798
$ws.Range("H15").Value = 0.1390134529147982
$ws.Range("I15").Value = 0.07174887892376682
$ws.Range("J15").Value = 0.3901345291479821
$ws.Range("K15").Value = 0.05381165919282511
$ws.Range("O15").Value = 0.04484304932735426
$ws.Range("S15").Value = 0.2690582959641256
$ws.Range("F16").Value = 0.01149425287356322
$ws.Range("H16").Value = 0.2011494252873563
$ws.Range("I16").Value = 0.06321839080459771
$ws.Range("J16").Value = 0.3908045977011494
$ws.Range("K16").Value = 0.09195402298850575
$ws.Range("M16").Value = 0.01149425287356322
$ws.Range("O16").Value = 0.07471264367816093
$ws.Range("S16").Value = 0.1551724137931035
$ws.Range("F17").Value = 0.02631578947368421
$ws.Range("H17").Value = 0.1842105263157895
$ws.Range("I17").Value = 0.06680161943319839
$ws.Range("J17").Value = 0.4190283400809717
$ws.Range("K17").Value = 0.08097165991902834
$ws.Range("M17").Value = 0.01417004048582996
$ws.Range("N17").Value = 0.002024291497975709
$ws.Range("O17").Value = 0.07489878542510121
$ws.Range("S17").Value = 0.131578947368421
$ws.Range("F18").Value = 0.02553191489361702
$ws.Range("H18").Value = 0.2085106382978723
$ws.Range("I18").Value = 0.09787234042553192
$ws.Range("J18").Value = 0.3787234042553191
$ws.Range("K18").Value = 0.1191489361702128
$ws.Range("M18").Value = 0.01276595744680851
$ws.Range("O18").Value = 0.03404255319148936
$ws.Range("S18").Value = 0.1234042553191489
$ws.Range("F19").Value = 0.02265625
$ws.Range("H19").Value = 0.2453125
$ws.Range("I19").Value = 0.07656250000000001
$ws.Range("J19").Value = 0.35703125
$ws.Range("K19").Value = 0.09531249999999999
$ws.Range("M19").Value = 0.02421875
$ws.Range("N19").Value = 0.00078125
$ws.Range("O19").Value = 0.06640625
$ws.Range("S19").Value = 0.11171875
